# Case_1_49 / lines_states.xlsx
# Commit: "contingencies with rene fine"
#
# Two new line entries ("line7", "line8") are inserted into the lines/states
# table right after "line6" (i.e. at sheet rows 8 and 9), pushing the
# existing "extr1".."extr8" rows down by two (to rows 10..17). A couple of
# the "in_service" flags among the extr rows also flip from FALSE to TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room: clone the formatting of the last existing data row (15)
# into the two brand-new rows (16 and 17) so they get the same cell style
# (border/alignment/bold index column) as the rest of the table. ---
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A15:E15").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)

# --- Final data for rows 8..17 (index, name, from_bus, to_bus, in_service) ---
$data = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "applied contingency line rows"
